# 第 8 次commit About 更新 excel 文件
#
# Adds a new "测试01" label at A1, and appends a second test block
# ("测试02" / "浏览器正常运行代码") below the existing browser-support
# table, re-listing 火狐 / 谷歌chrome / iPad safari / 微软 Edge with a
# merged "YES" result cell for each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: new heading "测试01" -------------------------------------------
$ws.Range("A1").Value = "测试01"

# --- Row 9: second section heading "测试02" -------------------------------
$ws.Range("A9").Value = "测试02"

# --- Row 10: merged sub-heading "浏览器正常运行代码" ----------------------
# Reuse the bold/center formatting already used for column-A headers
# (A3:A7), then drop the vertical-center so it matches a plain
# horizontal-center style, same as Excel's own "centered" cell style.
$ws.Range("A3").Copy()
$ws.Range("B10:C10").PasteSpecial(-4122)
$ws.Range("B10:C10").VerticalAlignment = -4107
$ws.Range("B10").Value = "浏览器正常运行代码"
$ws.Range("B10:C10").Merge()

# --- Rows 11-14: repeat the browser list with merged "YES" results -------
$ws.Range("A3").Copy()
$ws.Range("A11:A14").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B11:C14").PasteSpecial(-4122)

$ws.Range("A11").Value = "火狐"
$ws.Range("A12").Value = "谷歌chrome"
$ws.Range("A13").Value = "iPad safari"
$ws.Range("A14").Value = "微软 Edge"

$ws.Range("B11").Value = "YES"
$ws.Range("B12").Value = "YES"
$ws.Range("B13").Value = "YES"
$ws.Range("B14").Value = "YES"

$ws.Range("B11:C11").Merge()
$ws.Range("B12:C12").Merge()
$ws.Range("B13:C13").Merge()
$ws.Range("B14:C14").Merge()

# --- View: scroll down and land the selection on the new last row --------
$ws.Range("B14:C14").Select()
